# Add a "Save" column (H) to the s_vals sheet:
#  - H1 header "Save", styled like the other header cells (copy format from G1)
#  - H2 value 0, matching the numeric style used by the other data cells

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from G1 onto H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data value under the new header.
$ws.Range("H2").Value = 0
